# -----------------------------------------------------------------------
# Adds a new "2022-Q3" quarter sheet (inserted right after "总计" and
# before "2022-Q2"), populates it with fund-holding detail data, and
# inserts a matching summary row at the top of the "总计" (totals) sheet.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet, positioned before "2022-Q2"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2SheetRef = $wb.Worksheets.Item("2022-Q2")

$newSheet = $wb.Worksheets.Add($q2SheetRef)
$newSheet.Name = "2022-Q3"

# NOTE: once the new sheet is inserted at the position previously held by
# "2022-Q2", the old $q2SheetRef variable becomes an alias for the *new*
# sheet (references appear to be position-based). Re-fetch "2022-Q2" by
# name so subsequent copies pull formatting/values from the real sheet.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Copy sheet-level formatting (outline / page setup / header style) from
# an existing quarter sheet so the new sheet matches the others.
$totalSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Header row (bold / centered / bordered, matching other quarter sheets)
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$q2Sheet.Range("B1").Copy()
for ($col = 2; $col -le 8; $col++) {
    $newSheet.Cells.Item(1, $col).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# Data rows for 2022-Q3 (index column A, then B..H)
$data = @(
    @(0,  "005669", "前海开源公用事业行业股票",                 "168.24", "90.44", "7.11", "11.9619", 4),
    @(1,  "001875", "前海开源沪港深优势精选灵活配置混合A",       "70.70",  "92.88", "8.44", "5.9671",  3),
    @(2,  "010717", "前海开源优质企业6个月持有期混合A",          "42.55",  "92.85", "8.45", "3.5955",  4),
    @(3,  "001837", "前海开源沪港深蓝筹精选灵活配置混合A",       "11.54",  "92.76", "8.45", "0.9751",  3),
    @(4,  "001874", "前海开源沪港深价值精选灵活配置混合",        "7.45",   "92.87", "8.15", "0.6072",  4),
    @(5,  "013270", "前海开源聚利一年持有混合A",                 "6.54",   "92.89", "8.28", "0.5415",  3),
    @(6,  "010452", "广发瑞福精选混合A",                         "11.39",  "83.41", "4.62", "0.5262",  1),
    @(7,  "010718", "前海开源优质企业6个月持有期混合C",          "5.77",   "92.85", "8.45", "0.4876",  4),
    @(8,  "012943", "广发稳睿六个月持有期混合A",                 "20.75",  "26.11", "1.68", "0.3486",  5),
    @(9,  "012944", "广发稳睿六个月持有期混合C",                 "19.14",  "26.11", "1.68", "0.3216",  5),
    @(10, "011635", "富国港股通策略精选混合A",                   "6.21",   "73.36", "3.78", "0.2347",  7),
    @(11, "011481", "广发瑞锦一年定开混合",                      "2.62",   "89.31", "6.99", "0.1831",  3),
    @(12, "005197", "工银瑞信沪港深精选灵活配置混合A",           "4.01",   "94.39", "3.66", "0.1468",  10),
    @(13, "002653", "泰康沪港深精选灵活配置混合",                "7.08",   "86.37", "1.64", "0.1161",  9),
    @(14, "010453", "广发瑞福精选混合C",                         "1.09",   "83.41", "4.62", "0.0504",  1),
    @(15, "005198", "工银瑞信沪港深精选灵活配置混合C",           "1.01",   "94.39", "3.66", "0.0370",  10),
    @(16, "013271", "前海开源聚利一年持有混合C",                 "0.39",   "92.89", "8.28", "0.0323",  3),
    @(17, "011871", "前海开源沪港深优势精选灵活配置混合C",       "0.32",   "92.88", "8.44", "0.0270",  3),
    @(18, "003580", "泰康沪港深价值优选灵活配置混合",            "1.26",   "89.69", "1.75", "0.0220",  9),
    @(19, "011636", "富国港股通策略精选混合C",                   "0.57",   "73.36", "3.78", "0.0215",  7),
    @(20, "006923", "前海开源沪港深非周期性行业股票A",           "0.26",   "86.78", "5.36", "0.0139",  9),
    @(21, "006924", "前海开源沪港深非周期性行业股票C",           "0.20",   "86.78", "5.36", "0.0107",  9),
    @(22, "012711", "前海开源沪港深蓝筹精选灵活配置混合C",       "0.08",   "92.76", "8.45", "0.0068",  3),
    @(23, "004532", "民生加银中证港股通高股息精选指数A",         "0.13",   "92.87", "4.35", "0.0057",  5),
    @(24, "004533", "民生加银中证港股通高股息精选指数C",         "0.08",   "92.87", "4.35", "0.0035",  5)
)

$row = 2
foreach ($r in $data) {
    # Column A: numeric row index, styled like the other quarter sheets
    $q2Sheet.Range("A2").Copy()
    $newSheet.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $newSheet.Cells.Item($row, 1).Value = $r[0]

    # Column B: fund code -- must stay text so leading zeros survive
    $newSheet.Cells.Item($row, 2).NumberFormat = "@"
    $newSheet.Cells.Item($row, 2).Value = $r[1]

    # Column C: fund name (plain text)
    $newSheet.Cells.Item($row, 3).Value = $r[2]

    # Columns D-G: stored as text (matches source formatting)
    $newSheet.Cells.Item($row, 4).NumberFormat = "@"
    $newSheet.Cells.Item($row, 4).Value = $r[3]
    $newSheet.Cells.Item($row, 5).NumberFormat = "@"
    $newSheet.Cells.Item($row, 5).Value = $r[4]
    $newSheet.Cells.Item($row, 6).NumberFormat = "@"
    $newSheet.Cells.Item($row, 6).Value = $r[5]
    $newSheet.Cells.Item($row, 7).NumberFormat = "@"
    $newSheet.Cells.Item($row, 7).Value = $r[6]

    # Column H: numeric rank
    $newSheet.Cells.Item($row, 8).Value = $r[7]

    $row++
}

$newSheet.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Insert the new 2022-Q3 summary row into the "总计" sheet
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert() | Out-Null

# Normalise formatting of the freshly-inserted row: index column (A)
# gets the bold/centered/bordered style used by the other rows, while
# B:D stay unstyled like the existing data rows.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$totalSheet.Range("A2").Value = 0

$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 25
$totalSheet.Range("D2").Value = 26.24

# Re-number the index column (A) for the remaining rows: 0..6
for ($i = 0; $i -le 5; $i++) {
    $totalSheet.Cells.Item($i + 3, 1).Value = $i + 1
}

$totalSheet.Range("A1").Select() | Out-Null
